# Updating the model for MM&MV
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2 to 97) forward by one day
$ws.Cells.Item(2, 1).Value = 46070.01041666666
$ws.Cells.Item(3, 1).Value = 46070.02083333334
$ws.Cells.Item(4, 1).Value = 46070.03125
$ws.Cells.Item(5, 1).Value = 46070.04166666666
$ws.Cells.Item(6, 1).Value = 46070.05208333334
$ws.Cells.Item(7, 1).Value = 46070.0625
$ws.Cells.Item(8, 1).Value = 46070.07291666666
$ws.Cells.Item(9, 1).Value = 46070.08333333334
$ws.Cells.Item(10, 1).Value = 46070.09375
$ws.Cells.Item(11, 1).Value = 46070.10416666666
$ws.Cells.Item(12, 1).Value = 46070.11458333334
$ws.Cells.Item(13, 1).Value = 46070.125
$ws.Cells.Item(14, 1).Value = 46070.13541666666
$ws.Cells.Item(15, 1).Value = 46070.14583333334
$ws.Cells.Item(16, 1).Value = 46070.15625
$ws.Cells.Item(17, 1).Value = 46070.16666666666
$ws.Cells.Item(18, 1).Value = 46070.17708333334
$ws.Cells.Item(19, 1).Value = 46070.1875
$ws.Cells.Item(20, 1).Value = 46070.19791666666
$ws.Cells.Item(21, 1).Value = 46070.20833333334
$ws.Cells.Item(22, 1).Value = 46070.21875
$ws.Cells.Item(23, 1).Value = 46070.22916666666
$ws.Cells.Item(24, 1).Value = 46070.23958333334
$ws.Cells.Item(25, 1).Value = 46070.25
$ws.Cells.Item(26, 1).Value = 46070.26041666666
$ws.Cells.Item(27, 1).Value = 46070.27083333334
$ws.Cells.Item(28, 1).Value = 46070.28125
$ws.Cells.Item(29, 1).Value = 46070.29166666666
$ws.Cells.Item(30, 1).Value = 46070.30208333334
$ws.Cells.Item(31, 1).Value = 46070.3125
$ws.Cells.Item(32, 1).Value = 46070.32291666666
$ws.Cells.Item(33, 1).Value = 46070.33333333334
$ws.Cells.Item(34, 1).Value = 46070.34375
$ws.Cells.Item(35, 1).Value = 46070.35416666666
$ws.Cells.Item(36, 1).Value = 46070.36458333334
$ws.Cells.Item(37, 1).Value = 46070.375
$ws.Cells.Item(38, 1).Value = 46070.38541666666
$ws.Cells.Item(39, 1).Value = 46070.39583333334
$ws.Cells.Item(40, 1).Value = 46070.40625
$ws.Cells.Item(41, 1).Value = 46070.41666666666
$ws.Cells.Item(42, 1).Value = 46070.42708333334
$ws.Cells.Item(43, 1).Value = 46070.4375
$ws.Cells.Item(44, 1).Value = 46070.44791666666
$ws.Cells.Item(45, 1).Value = 46070.45833333334
$ws.Cells.Item(46, 1).Value = 46070.46875
$ws.Cells.Item(47, 1).Value = 46070.47916666666
$ws.Cells.Item(48, 1).Value = 46070.48958333334
$ws.Cells.Item(49, 1).Value = 46070.5
$ws.Cells.Item(50, 1).Value = 46070.51041666666
$ws.Cells.Item(51, 1).Value = 46070.52083333334
$ws.Cells.Item(52, 1).Value = 46070.53125
$ws.Cells.Item(53, 1).Value = 46070.54166666666
$ws.Cells.Item(54, 1).Value = 46070.55208333334
$ws.Cells.Item(55, 1).Value = 46070.5625
$ws.Cells.Item(56, 1).Value = 46070.57291666666
$ws.Cells.Item(57, 1).Value = 46070.58333333334
$ws.Cells.Item(58, 1).Value = 46070.59375
$ws.Cells.Item(59, 1).Value = 46070.60416666666
$ws.Cells.Item(60, 1).Value = 46070.61458333334
$ws.Cells.Item(61, 1).Value = 46070.625
$ws.Cells.Item(62, 1).Value = 46070.63541666666
$ws.Cells.Item(63, 1).Value = 46070.64583333334
$ws.Cells.Item(64, 1).Value = 46070.65625
$ws.Cells.Item(65, 1).Value = 46070.66666666666
$ws.Cells.Item(66, 1).Value = 46070.67708333334
$ws.Cells.Item(67, 1).Value = 46070.6875
$ws.Cells.Item(68, 1).Value = 46070.69791666666
$ws.Cells.Item(69, 1).Value = 46070.70833333334
$ws.Cells.Item(70, 1).Value = 46070.71875
$ws.Cells.Item(71, 1).Value = 46070.72916666666
$ws.Cells.Item(72, 1).Value = 46070.73958333334
$ws.Cells.Item(73, 1).Value = 46070.75
$ws.Cells.Item(74, 1).Value = 46070.76041666666
$ws.Cells.Item(75, 1).Value = 46070.77083333334
$ws.Cells.Item(76, 1).Value = 46070.78125
$ws.Cells.Item(77, 1).Value = 46070.79166666666
$ws.Cells.Item(78, 1).Value = 46070.80208333334
$ws.Cells.Item(79, 1).Value = 46070.8125
$ws.Cells.Item(80, 1).Value = 46070.82291666666
$ws.Cells.Item(81, 1).Value = 46070.83333333334
$ws.Cells.Item(82, 1).Value = 46070.84375
$ws.Cells.Item(83, 1).Value = 46070.85416666666
$ws.Cells.Item(84, 1).Value = 46070.86458333334
$ws.Cells.Item(85, 1).Value = 46070.875
$ws.Cells.Item(86, 1).Value = 46070.88541666666
$ws.Cells.Item(87, 1).Value = 46070.89583333334
$ws.Cells.Item(88, 1).Value = 46070.90625
$ws.Cells.Item(89, 1).Value = 46070.91666666666
$ws.Cells.Item(90, 1).Value = 46070.92708333334
$ws.Cells.Item(91, 1).Value = 46070.9375
$ws.Cells.Item(92, 1).Value = 46070.94791666666
$ws.Cells.Item(93, 1).Value = 46070.95833333334
$ws.Cells.Item(94, 1).Value = 46070.96875
$ws.Cells.Item(95, 1).Value = 46070.97916666666
$ws.Cells.Item(96, 1).Value = 46070.98958333334
$ws.Cells.Item(97, 1).Value = 46071

# Update the three changed production values in column B
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(29, 2).Value = 11
$ws.Cells.Item(30, 2).Value = 27
